# Update PSSM data values on Sheet1 (supplemental figures update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -17.56123572887572
$ws.Range("C2").Value = 1.926033182688413
$ws.Range("D2").Value = -17.56123572887572
$ws.Range("E2").Value = -17.56123572887572
$ws.Range("F2").Value = -17.56123572887572
$ws.Range("G2").Value = -17.56123572887572
$ws.Range("H2").Value = -17.56123572887572
$ws.Range("I2").Value = -17.56123572887572
$ws.Range("J2").Value = -17.56123572887572
$ws.Range("K2").Value = -17.56123572887572

# Row 3
$ws.Range("B3").Value = -17.56123572887572
$ws.Range("C3").Value = -17.56123572887572
$ws.Range("D3").Value = -17.56123572887572
$ws.Range("E3").Value = -17.56123572887572
$ws.Range("F3").Value = -17.56123572887572
$ws.Range("G3").Value = -17.56123572887572
$ws.Range("H3").Value = -17.56123572887572
$ws.Range("I3").Value = 2.545916202021252
$ws.Range("J3").Value = -17.56123572887572
$ws.Range("K3").Value = -17.56123572887572

# Row 4
$ws.Range("B4").Value = -17.56123572887572
$ws.Range("C4").Value = 2.192289537341488
$ws.Range("D4").Value = 2.099848197142211
$ws.Range("E4").Value = -17.56123572887572
$ws.Range("F4").Value = 3.49177335923862
$ws.Range("G4").Value = -17.56123572887572
$ws.Range("H4").Value = 1.973266644299622
$ws.Range("I4").Value = -17.56123572887572
$ws.Range("J4").Value = 2.21470403027035
$ws.Range("K4").Value = -17.56123572887572

# Row 5
$ws.Range("B5").Value = -17.56123572887572
$ws.Range("C5").Value = 1.915037804282477
$ws.Range("D5").Value = -17.56123572887572
$ws.Range("E5").Value = -17.56123572887572
$ws.Range("F5").Value = -17.56123572887572
$ws.Range("G5").Value = 2.787692291224349
$ws.Range("H5").Value = -17.56123572887572
$ws.Range("I5").Value = -17.56123572887572
$ws.Range("J5").Value = -17.56123572887572
$ws.Range("K5").Value = -17.56123572887572

# Row 6
$ws.Range("B6").Value = -17.56123572887572
$ws.Range("C6").Value = -17.56123572887572
$ws.Range("D6").Value = -17.56123572887572
$ws.Range("E6").Value = -17.56123572887572
$ws.Range("F6").Value = -17.56123572887572
$ws.Range("G6").Value = -17.56123572887572
$ws.Range("H6").Value = -17.56123572887572
$ws.Range("I6").Value = -17.56123572887572
$ws.Range("J6").Value = -17.56123572887572
$ws.Range("K6").Value = -17.56123572887572

# Row 7
$ws.Range("B7").Value = 2.591406860690374
$ws.Range("C7").Value = -17.56123572887572
$ws.Range("D7").Value = -17.56123572887572
$ws.Range("E7").Value = -17.56123572887572
$ws.Range("F7").Value = -17.56123572887572
$ws.Range("G7").Value = -17.56123572887572
$ws.Range("H7").Value = -17.56123572887572
$ws.Range("I7").Value = -17.56123572887572
$ws.Range("J7").Value = -17.56123572887572
$ws.Range("K7").Value = -17.56123572887572

# Row 8
$ws.Range("B8").Value = -17.56123572887572
$ws.Range("C8").Value = -17.56123572887572
$ws.Range("D8").Value = -17.56123572887572
$ws.Range("E8").Value = 1.782672142466858
$ws.Range("F8").Value = -17.56123572887572
$ws.Range("G8").Value = -17.56123572887572
$ws.Range("H8").Value = -17.56123572887572
$ws.Range("I8").Value = -17.56123572887572
$ws.Range("J8").Value = -17.56123572887572
$ws.Range("K8").Value = -17.56123572887572

# Row 9
$ws.Range("B9").Value = 3.80457460755725
$ws.Range("C9").Value = -17.56123572887572
$ws.Range("D9").Value = -17.56123572887572
$ws.Range("E9").Value = -17.56123572887572
$ws.Range("F9").Value = -17.56123572887572
$ws.Range("G9").Value = -17.56123572887572
$ws.Range("H9").Value = -17.56123572887572
$ws.Range("I9").Value = -17.56123572887572
$ws.Range("J9").Value = -17.56123572887572
$ws.Range("K9").Value = -17.56123572887572

# Row 10
$ws.Range("B10").Value = -17.56123572887572
$ws.Range("C10").Value = -17.56123572887572
$ws.Range("D10").Value = -17.56123572887572
$ws.Range("E10").Value = -17.56123572887572
$ws.Range("F10").Value = -17.56123572887572
$ws.Range("G10").Value = -17.56123572887572
$ws.Range("H10").Value = -17.56123572887572
$ws.Range("I10").Value = 1.309396894222689
$ws.Range("J10").Value = -17.56123572887572
$ws.Range("K10").Value = -17.56123572887572

# Row 11
$ws.Range("B11").Value = -17.56123572887572
$ws.Range("C11").Value = -17.56123572887572
$ws.Range("D11").Value = -17.56123572887572
$ws.Range("E11").Value = 2.884794488641071
$ws.Range("F11").Value = -17.56123572887572
$ws.Range("G11").Value = 2.888163832866709
$ws.Range("H11").Value = -17.56123572887572
$ws.Range("I11").Value = -17.56123572887572
$ws.Range("J11").Value = -17.56123572887572
$ws.Range("K11").Value = -17.56123572887572

# Row 12
$ws.Range("B12").Value = -17.56123572887572
$ws.Range("C12").Value = -17.56123572887572
$ws.Range("D12").Value = -17.56123572887572
$ws.Range("E12").Value = -17.56123572887572
$ws.Range("F12").Value = -17.56123572887572
$ws.Range("G12").Value = -17.56123572887572
$ws.Range("H12").Value = -17.56123572887572
$ws.Range("I12").Value = -17.56123572887572
$ws.Range("J12").Value = -17.56123572887572
$ws.Range("K12").Value = -17.56123572887572

# Row 13
$ws.Range("B13").Value = -17.56123572887572
$ws.Range("C13").Value = -17.56123572887572
$ws.Range("D13").Value = -17.56123572887572
$ws.Range("E13").Value = 2.399469058353237
$ws.Range("F13").Value = -17.56123572887572
$ws.Range("G13").Value = -17.56123572887572
$ws.Range("H13").Value = -17.56123572887572
$ws.Range("I13").Value = -17.56123572887572
$ws.Range("J13").Value = 1.905073784292669
$ws.Range("K13").Value = -17.56123572887572

# Row 14
$ws.Range("B14").Value = -17.56123572887572
$ws.Range("C14").Value = -17.56123572887572
$ws.Range("D14").Value = 1.178020780522898
$ws.Range("E14").Value = -17.56123572887572
$ws.Range("F14").Value = -17.56123572887572
$ws.Range("G14").Value = -17.56123572887572
$ws.Range("H14").Value = -17.56123572887572
$ws.Range("I14").Value = -17.56123572887572
$ws.Range("J14").Value = -17.56123572887572
$ws.Range("K14").Value = 4.321921008245333

# Row 15
$ws.Range("B15").Value = -17.56123572887572
$ws.Range("C15").Value = -17.56123572887572
$ws.Range("D15").Value = 1.080015200721807
$ws.Range("E15").Value = -17.56123572887572
$ws.Range("F15").Value = -17.56123572887572
$ws.Range("G15").Value = -17.56123572887572
$ws.Range("H15").Value = -17.56123572887572
$ws.Range("I15").Value = -17.56123572887572
$ws.Range("J15").Value = -17.56123572887572
$ws.Range("K15").Value = -17.56123572887572

# Row 16
$ws.Range("B16").Value = -17.56123572887572
$ws.Range("C16").Value = -17.56123572887572
$ws.Range("D16").Value = -17.56123572887572
$ws.Range("E16").Value = -17.56123572887572
$ws.Range("F16").Value = -17.56123572887572
$ws.Range("G16").Value = -17.56123572887572
$ws.Range("H16").Value = -17.56123572887572
$ws.Range("I16").Value = -17.56123572887572
$ws.Range("J16").Value = 2.225463041310161
$ws.Range("K16").Value = -17.56123572887572

# Row 17
$ws.Range("B17").Value = -17.56123572887572
$ws.Range("C17").Value = 2.015870347971362
$ws.Range("D17").Value = 2.441520238816671
$ws.Range("E17").Value = -17.56123572887572
$ws.Range("F17").Value = -17.56123572887572
$ws.Range("G17").Value = -17.56123572887572
$ws.Range("H17").Value = 1.213837581975989
$ws.Range("I17").Value = 2.043999610560206
$ws.Range("J17").Value = 2.220637814819483
$ws.Range("K17").Value = -17.56123572887572

# Row 18
$ws.Range("B18").Value = -17.56123572887572
$ws.Range("C18").Value = -17.56123572887572
$ws.Range("D18").Value = -17.56123572887572
$ws.Range("E18").Value = -17.56123572887572
$ws.Range("F18").Value = -17.56123572887572
$ws.Range("G18").Value = -17.56123572887572
$ws.Range("H18").Value = 1.254863557829424
$ws.Range("I18").Value = 1.352392934718297
$ws.Range("J18").Value = 1.185977123567607
$ws.Range("K18").Value = -17.56123572887572

# Row 19
$ws.Range("B19").Value = -17.56123572887572
$ws.Range("C19").Value = -17.56123572887572
$ws.Range("D19").Value = 1.591606120867319
$ws.Range("E19").Value = -17.56123572887572
$ws.Range("F19").Value = -17.56123572887572
$ws.Range("G19").Value = -17.56123572887572
$ws.Range("H19").Value = 1.613175546917496
$ws.Range("I19").Value = 1.524322311019658
$ws.Range("J19").Value = -17.56123572887572
$ws.Range("K19").Value = -17.56123572887572

# Row 20
$ws.Range("B20").Value = -17.56123572887572
$ws.Range("C20").Value = 0.7942917358732815
$ws.Range("D20").Value = 1.531358839591237
$ws.Range("E20").Value = -17.56123572887572
$ws.Range("F20").Value = 3.129370185122707
$ws.Range("G20").Value = -17.56123572887572
$ws.Range("H20").Value = 1.94044792831553
$ws.Range("I20").Value = 1.089511022806619
$ws.Range("J20").Value = -17.56123572887572
$ws.Range("K20").Value = -17.56123572887572

# Row 21
$ws.Range("B21").Value = -17.56123572887572
$ws.Range("C21").Value = 1.056464887080042
$ws.Range("D21").Value = -17.56123572887572
$ws.Range("E21").Value = 1.962435334831125
$ws.Range("F21").Value = -17.56123572887572
$ws.Range("G21").Value = 2.5087816754502
$ws.Range("H21").Value = 2.160239001496894
$ws.Range("I21").Value = -17.56123572887572
$ws.Range("J21").Value = -17.56123572887572
$ws.Range("K21").Value = -17.56123572887572

